$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Replace the two "sample description N" placeholder rows with the
#    real listing descriptions (house/agent primary-key clean-up).
#    D4 is written before D3 so the shared-string table picks up the
#    "2 Bedroom 1 Bath..." text ahead of the "Welcome home..." text,
#    matching the canonical save order.
# ------------------------------------------------------------------

$d4Text = @'
2 Bedroom 1 Bath complete remodel. Updated kitchen: center island, New cabinets, New appliances, Granite counter tops. Great views of the backyard from the kitchen and breakfast room windows! Remodel and new appliances 2023. Open concept. Vaulted Living area. Tall baseboards through out, and New Luxury Vinyl plank flooring. New windows to help keep your utility bills down, hot water heater. HVAC recently services and new AC outside unit (condenser) installed. Large backyard. Close to Six Flags, Hurricane Harbor, Globe Life Field, AT&T stadium, and more! Great Location!
'@

$ws.Range("D4").Value = $d4Text
$ws.Range("D4").WrapText = $true
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").VerticalAlignment = -4108
$ws.Rows(4).RowHeight = 144

$d3Text = @'
Welcome home to this brand new bi-level condo on Aldine! This condo is perfectly positioned to access all of the great restaurants and retail shops in the area! Enter unit A to a stunning open-concept living area that leads into a wrap-around kitchen equipped with quartz counters, stainless steel appliances and marble backsplash. Bedrooms on both floors offer tons of space and have direct access to a private bath.

'@

$ws.Range("D3").Value = $d3Text
$ws.Range("D3").WrapText = $true
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4108
$ws.Rows(3).RowHeight = 115.2

# ------------------------------------------------------------------
# 2) Update selection / scroll position of the sheet view.
# ------------------------------------------------------------------
$ws.Range("G4").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

# ------------------------------------------------------------------
# 3) Add the "Links for descriptions" textbox/shape that documents
#    where each listing description came from.
# ------------------------------------------------------------------
$shp = $ws.Shapes.AddTextbox(1, 50.4, 552, 603, 214.8)
$shp.Name = "TextBox 1"

$shp.Fill.ForeColor.RGB = RGB(255, 255, 255)
$shp.Line.ForeColor.RGB = RGB(127, 127, 127)
$shp.Line.Weight = 0.75

$tb = $shp.TextFrame
$tb.Characters.Text = "Links for descriptions`r`n`r`n1. 4402 Driftwood Dr, Philadelphia, PA 19129 | Zillow`r`n2. 407 W Oxford St, Philadelphia, PA 19122, USA - 3 unit Rentals | Zumper`r`n3. 2501 Oak Hill Drive, Arlington, TX 76006, USA | 2 bed duplex for rent #90051616 | Rentberry"

# ------------------------------------------------------------------
# 4) Attach the hyperlinks for each of the three listed sources.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($shp, "https://www.zillow.com/homedetails/4402-Driftwood-Dr-Philadelphia-PA-19129", "", "", "4402 Driftwood Dr, Philadelphia, PA 19129 | Zillow")
$ws.Hyperlinks.Add($shp, "https://www.zumper.com/apartments-for-rent/407-w-oxford-st-philadelphia-pa-19122-usa", "", "", "407 W Oxford St, Philadelphia, PA 19122, USA - 3 unit Rentals | Zumper")
$ws.Hyperlinks.Add($shp, "https://rentberry.com/listings/90051616-2501-oak-hill-drive-arlington-tx-76006-usa", "", "", "2501 Oak Hill Drive, Arlington, TX 76006, USA | 2 bed duplex for rent #90051616 | Rentberry")
